$d = $word.ActiveDocument

# Locate the "GitHub Desktop : " paragraph (last paragraph of the body) that
# needs to be reformatted, split into new runs, and followed by several new
# paragraphs describing Visual Studio Code and Perplexity.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("GitHub Desktop")) {
        $target = $p
    }
}

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:spacing w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Git</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Hub Desktop</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t> </w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>C’est une</w:t></w:r><w:r><w:t> </w:t></w:r><w:r><w:t xml:space="preserve">interface graphique pour Git permettant le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>versionnement</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, que ce soit pour les différents documents ou le code source. On peut également</w:t></w:r><w:r><w:t xml:space="preserve"> sélectionner les versions précédentes </w:t></w:r><w:r><w:t>en cas de bug ou autres problèmes.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:spacing w:val="1"/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Visual studio code</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:b/><w:bCs/><w:spacing w:val="1"/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:spacing w:val="1"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>c’est un éditeur de code qui nous permet de créer le code pour le projet INFOPROD. Il dispose de</w:t></w:r><w:r><w:t xml:space="preserve"> plusieurs extensions intégrées </w:t></w:r><w:r><w:t xml:space="preserve">(PHP, </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>MySQL,</w:t></w:r><w:r><w:t>GIT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>) pour le développement web et la gestion de base de données.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Perplexity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t> :</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>C’est</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>intelligence artificielle qui m’aide, que ce soit pour corriger des bugs ou simplement comprendre des sujets autour du projet INFOPROD</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:b/><w:bCs/><w:spacing w:val="1"/></w:rPr></w:pPr></w:p><w:p/>
'@

$target.Range.InsertXML($xml)
